$d = $word.ActiveDocument

# 1. Street address: "VALENTIN GOMEZ FARIAS # 208 BARRIO SAN LORENZO" -> "VALENTIN GOMEZ FARIAS #2031"
$d.Content.Find.Execute("VALENTIN GOMEZ FARIAS # 208 BARRIO SAN LORENZO", $true, $false, $false, $false, $false, $true, 1, $false, "VALENTIN GOMEZ FARIAS #2031", 2)

# 2. City line: " ZIMATLAN DE ALVAREZ, OAX." -> " BARRIO SAN LORENZO, ZIMATLÁN DE ÁLVAREZ, OAXACA."
$d.Content.Find.Execute(" ZIMATLAN DE ALVAREZ, OAX.", $true, $false, $false, $false, $false, $true, 1, $false, " BARRIO SAN LORENZO, ZIMATLÁN DE ÁLVAREZ, OAXACA.", 2)

# 3. Postal code "C. P." value "0" -> "71200"
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("C. P.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cpValRange = $d.Range($rng.End + 1, $rng.End + 2)
$cpValRange.Text = "71200"

# 4. Date change
$d.Content.Find.Execute("20 DE JUNIO DE 2025", $true, $false, $false, $false, $false, $true, 1, $false, "14 DE NOVIEMBRE DE 2025", 2)
